$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-11: Colaborador_id, Colaborador_nome, Departamento, Motivo_da_ausencia,
# Horas_de_ausencia, Data_da_ausencia (serial date), Salario
$data = @(
    @(2,  58228, "Bernardo Monteiro",    "TI",                      "Viagem de negócios", 4, 45097, 5687.87),
    @(3,  43536, "Davi Lucas Rodrigues", "P&D",                     "Outros",             4, 45090, 5706.39),
    @(4,  30574, "Benício da Luz",       "Vendas",                  "Doença",             7, 45081, 10893.35),
    @(5,  77705, "Luana Novaes",         "Atendimento ao Cliente",  "Consulta médica",    2, 45097, 11911.82),
    @(6,  33035, "Ana Júlia Caldeira",   "Vendas",                  "Outros",             6, 45092, 8607.47),
    @(7,  52175, "Mirella Azevedo",      "Recursos Humanos",        "Consulta médica",    6, 45078, 7361.52),
    @(8,  33275, "Ana Gomes",            "Jurídico",                "Problemas pessoais", 8, 45102, 9415.29),
    @(9,  77513, "Vitor Hugo Caldeira",  "Jurídico",                "Problemas pessoais", 4, 45086, 8075.54),
    @(10, 99008, "Nathan Souza",         "TI",                      "Problemas pessoais", 2, 45098, 7522.65),
    @(11, 69934, "João da Conceição",    "Jurídico",                "Viagem de negócios", 4, 45100, 2595.5)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}
